$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Narrow column B (47 -> 43 characters)
$ws.Columns("B").ColumnWidth = 42.14

# Drop every existing hyperlink; fresh ones get re-added below for the rows that remain
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-09-10 06:25:20'
$ws.Range("B2").Value = '健康分野でのAIローカルサーバー構築'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5390534'
$ws.Range("G2").Value = 303
$ws.Range("H2").Value = '🔥AI,Ai'

# Row 3
$ws.Range("A3").Value = '2025-09-10 06:25:20'
$ws.Range("B3").Value = '【GAS開発】Meta広告・YouTubeインサイトデータ取得&動画自動投稿ツール'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5390748'
$ws.Range("G3").Value = 163
$ws.Range("H3").Value = '◆ツール,開発 ◇サイト'

# Row 4
$ws.Range("A4").Value = '2025-09-10 06:25:20'
$ws.Range("B4").Value = '【n8n等最先端ツールを使った】業務自動化ブログのTech記事ライター'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5390712'
$ws.Range("G4").Value = 140
$ws.Range("H4").Value = '◆ツール,自動化'

# Row 5
$ws.Range("A5").Value = '2025-09-10 06:25:20'
$ws.Range("B5").Value = 'システムの開発補助や運営サポート【フルリモート×長期】'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5323359'
$ws.Range("G5").Value = 83
$ws.Range("H5").Value = '◆開発'

# Row 6
$ws.Range("A6").Value = '2025-09-10 06:25:20'
$ws.Range("B6").Value = '【急募】ファクタリング会社の会員ページ開発を依頼します'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5390814'
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = '◆開発'

# Row 7
$ws.Range("A7").Value = '2025-09-10 06:25:20'
$ws.Range("B7").Value = '【急募】自社アプリのデバッグ・バグチェック業務依頼 ※NDA締結必須'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '~ 5,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5390852'
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = '◇アプリ'

# Row 8
$ws.Range("A8").Value = '2025-09-10 06:25:20'
$ws.Range("B8").Value = '【フリーランス歓迎】フロントエンドエンジニア募集(長期参画/リモート併用)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5390888'
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# Row 9
$ws.Range("A9").Value = '2025-09-10 06:25:20'
$ws.Range("B9").Value = '急募 限定公開 限定公開の仕事'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5390577'
$ws.Range("G9").Value = 25
$ws.Range("H9").ClearContents()

# Row 10
$ws.Range("A10").Value = '2025-09-10 06:25:20'
$ws.Range("B10").Value = '限定公開 PR 限定公開の仕事'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5385681'
$ws.Range("G10").Value = 25
$ws.Range("H10").ClearContents()

# Row 11
$ws.Range("A11").Value = '2025-09-10 06:25:20'
$ws.Range("B11").Value = '【急募】iPhone SwiftUIで魅力的なアニメーション制作依頼'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5390603'
$ws.Range("G11").Value = 18
$ws.Range("H11").ClearContents()

# Row 12
$ws.Range("A12").Value = '2025-09-10 06:25:20'
$ws.Range("B12").Value = '【急募】eLTAX導入サポート'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5390906'
$ws.Range("G12").Value = 10
$ws.Range("H12").ClearContents()

# Drop the now-stale rows 13-18 (list shrinks from 18 to 12 data rows -> dimension A1:H12)
$ws.Rows("13:18").Delete()

# Re-create hyperlinks on the URL column for the surviving rows, preserving the Hyperlink style
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5390534')
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5390748')
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5390712')
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5323359')
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5390814')
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5390852')
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5390888')
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5390577')
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5385681')
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5390603')
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5390906')
$ws.Range("F12").Style = "Hyperlink"
